$wb = $excel.ActiveWorkbook

# --- 2_UserManagement sheet: update the two data rows with new test-user info ---
$ws3 = $wb.Worksheets.Item("2_UserManagement")

# Row 2 (write order matches the new shared-string creation order: Tenant, checkbox, email)
$ws3.Range("K2").Value = "SVB Founders Platform"
$ws3.Range("J2").Value = "'true"
$ws3.Range("H2").Value = "xxx@janeirodigital.com"
$ws3.Range("I2").Value = "xxx@janeirodigital.com"

# Row 3
$ws3.Range("K3").Value = "SVB Founders Platform"
$ws3.Range("J3").Value = "'true"
$ws3.Range("H3").Value = "zzz@janeirodigital.com"
$ws3.Range("I3").Value = "zzz@janeirodigital.com"

# Move the view/selection on this sheet and make it the active (tab-selected) sheet
$ws3.Activate()
$window = $excel.ActiveWindow
$window.ScrollColumn = 8
$window.ScrollRow = 1
$ws3.Range("H3").Select()
